$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update "can exchange" flags (G column) from 0 -> 1 ---
$ws.Range("G23").Value2 = 1
$ws.Range("G29").Value2 = 1
$ws.Range("G30").Value2 = 1

# --- Fill in the new row 36 (2025 - 550th Anniversary - Birth of Michelangelo) ---
$ws.Range("B36").Value2 = "550th Anniversary - Birth of Michelangelo"

# D36 needs both its text value and the banding format copied from D35
$ws.Range("D35").Copy()
$ws.Range("D36").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D36").Value2 = "Obv: With mint symbol - ""R"""

$ws.Range("E36").Value2 = "Rev: new map of Europe"

# F36 must stay text (matching the existing "59.000" mintage label already used in F33)
# without Excel auto-converting it to a number, so copy the value only from F33.
$ws.Range("F33").Copy()
$ws.Range("F36").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("G36").Value2 = 1

$excel.CutCopyMode = $false

# --- Update sheet selection to B13 ---
$ws.Range("B13").Select()
